# chore: update Sheets via scheduled runner
# Refresh currentAveragePrice / LevePrice / LeveProfit figures on each
# item leve sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect the
# latest market-board pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 269.14
$ws.Range("I15").Value = 269.14
$ws.Range("K15").Value = 807.42
$ws.Range("M15").Value = -638.42

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2678.9092
$ws.Range("I138").Value = 1117.5333
$ws.Range("J138").Value = 2957.726
$ws.Range("K138").Value = 3352.5999
$ws.Range("L138").Value = 8873.178
$ws.Range("M138").Value = 1787.4001
$ws.Range("N138").Value = -19153.178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6260381
$ws.Range("I32").Value = 7151642
$ws.Range("J32").Value = 21555.3
$ws.Range("K32").Value = 7151642
$ws.Range("L32").Value = 21555.3
$ws.Range("M32").Value = -7151355
$ws.Range("N32").Value = -22129.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 82414.91
$ws.Range("I110").Value = 129087.43
$ws.Range("J110").Value = 738
$ws.Range("K110").Value = 129087.43
$ws.Range("L110").Value = 738
$ws.Range("M110").Value = -127042.43
$ws.Range("N110").Value = -4828

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 549.375
$ws.Range("I22").Value = 577.8570999999999
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 577.8570999999999
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -404.8570999999999
$ws.Range("N22").Value = -696

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 41668170
$ws.Range("I105").Value = 62501250
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 62501250
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -62499503
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 72709.14
$ws.Range("I107").Value = 101162.9
$ws.Range("J107").Value = 1574.75
$ws.Range("K107").Value = 101162.9
$ws.Range("L107").Value = 1574.75
$ws.Range("M107").Value = -99242.89999999999
$ws.Range("N107").Value = -5414.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1668.125
$ws.Range("I16").Value = 1686.5
$ws.Range("J16").Value = 1649.75
$ws.Range("K16").Value = 1686.5
$ws.Range("L16").Value = 1649.75
$ws.Range("M16").Value = -1399.5
$ws.Range("N16").Value = -2223.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5352.6914
$ws.Range("I31").Value = 1387.9412
$ws.Range("J31").Value = 8220.808999999999
$ws.Range("K31").Value = 1387.9412
$ws.Range("L31").Value = 8220.808999999999
$ws.Range("M31").Value = -1092.9412
$ws.Range("N31").Value = -8810.808999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5352.6914
$ws.Range("I34").Value = 1387.9412
$ws.Range("J34").Value = 8220.808999999999
$ws.Range("K34").Value = 1387.9412
$ws.Range("L34").Value = 8220.808999999999
$ws.Range("M34").Value = -1185.9412
$ws.Range("N34").Value = -8624.808999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1668.125
$ws.Range("I113").Value = 1686.5
$ws.Range("J113").Value = 1649.75
$ws.Range("K113").Value = 1686.5
$ws.Range("L113").Value = 1649.75
$ws.Range("M113").Value = 483.5
$ws.Range("N113").Value = -5989.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 281.1905
$ws.Range("I6").Value = 69.84614999999999
$ws.Range("J6").Value = 624.625
$ws.Range("K6").Value = 209.53845
$ws.Range("L6").Value = 1873.875
$ws.Range("M6").Value = -96.53844999999998
$ws.Range("N6").Value = -2099.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 4591
$ws.Range("J16").Value = 8882
$ws.Range("L16").Value = 26646
$ws.Range("N16").Value = -26992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1265.8334
$ws.Range("I57").Value = 500
$ws.Range("J57").Value = 1299.1305
$ws.Range("K57").Value = 1500
$ws.Range("L57").Value = 3897.3915
$ws.Range("M57").Value = -941
$ws.Range("N57").Value = -5015.3915

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 630.5
$ws.Range("I92").Value = 440
$ws.Range("J92").Value = 725.75
$ws.Range("K92").Value = 1320
$ws.Range("L92").Value = 2177.25
$ws.Range("M92").Value = -72
$ws.Range("N92").Value = -4673.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 325.77777
$ws.Range("I107").Value = 322.75
$ws.Range("K107").Value = 322.75
$ws.Range("M107").Value = 1597.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 66273.17999999999
$ws.Range("I113").Value = 70227.75
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 70227.75
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -68057.75
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 40000
$ws.Range("J50").Value = 40000
$ws.Range("L50").Value = 40000
$ws.Range("N50").Value = -41274

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 25002674
$ws.Range("I82").Value = 62502930
$ws.Range("J82").Value = 2504.5
$ws.Range("K82").Value = 62502930
$ws.Range("L82").Value = 2504.5
$ws.Range("M82").Value = -62502569
$ws.Range("N82").Value = -3226.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 25002674
$ws.Range("I85").Value = 62502930
$ws.Range("J85").Value = 2504.5
$ws.Range("K85").Value = 62502930
$ws.Range("L85").Value = 2504.5
$ws.Range("M85").Value = -62501682
$ws.Range("N85").Value = -5000.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 149111
$ws.Range("I62").Value = 204799.8
$ws.Range("K62").Value = 204799.8
$ws.Range("M62").Value = -204175.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 149111
$ws.Range("I65").Value = 204799.8
$ws.Range("K65").Value = 1023999
$ws.Range("M65").Value = -1020879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3698.0625
$ws.Range("I81").Value = 3530.8333
$ws.Range("J81").Value = 4199.75
$ws.Range("K81").Value = 7061.6666
$ws.Range("L81").Value = 8399.5
$ws.Range("M81").Value = -6000.6666
$ws.Range("N81").Value = -10521.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3698.0625
$ws.Range("I84").Value = 3530.8333
$ws.Range("J84").Value = 4199.75
$ws.Range("K84").Value = 35308.333
$ws.Range("L84").Value = 41997.5
$ws.Range("M84").Value = -30004.333
$ws.Range("N84").Value = -52605.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 4020000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2798.0908
$ws.Range("I136").Value = 2538.5715
$ws.Range("J136").Value = 3252.25
$ws.Range("K136").Value = 7615.7145
$ws.Range("L136").Value = 9756.75
$ws.Range("M136").Value = -5065.7145
$ws.Range("N136").Value = -14856.75
